$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$packages = @("autoawq", "duckdb", "pyjwt", "dash-extensions", "io", "getpadd", "jwcrypto")
$users = @("A", "B", "C", "D", "E", "F", "G")

$ws.Range("A1").Value = "Package"
$ws.Range("B1").Value = "User"

for ($i = 0; $i -lt $packages.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $packages[$i]
    $ws.Cells.Item($row, 2).Value = $users[$i]
}

$ws.Range("A2:A8").HorizontalAlignment = -4131

$excel.StandardFontSize = 11

$ws.Range("D8").Select()
